$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set cell values for new rows 56-61 ---
# Row 56 values
$ws.Cells.Item(56,1).Value = 43543.66878363426
$ws.Cells.Item(56,2).Value = "Entre 2 et 5 ans"
$ws.Cells.Item(56,3).Value = "Conférences (Jancovici, Bihouix etc.), Articles de vulgarisation & blogs, Livres, Vidéos Youtube de vulgarisation, Articles de presse, Publications scientifiques"
$ws.Cells.Item(56,4).Value = 2.0
$ws.Cells.Item(56,5).Value = 5.0
$ws.Cells.Item(56,6).Value = "Va commencer dans les 15 à 20 ans qui viennent"
$ws.Cells.Item(56,7).Value = "Un peu plus lent (de l'ordre de 20 à 30 ans)"
$ws.Cells.Item(56,8).Value = "Angoisse"
$ws.Cells.Item(56,9).Value = "On a un fort potentiel d'action à l'échelle individuelle, Je suis prêt à baisser mon niveau de vie si cette baisse s'opère pour les autres également, Je suis prêt à baisser mon niveau de vie même si cette baisse ne s'opère pas pour les autres, Une transition écologique efficace peut se faire dans un cadre démocratique, Pour agir efficacement, il faut hiérarchiser les problèmes (perte de biodiversité, dérèglement climatique etc.)"
$ws.Cells.Item(56,10).Value = "Une réaction de personnes ne voulant pas toucher à leur mode de vie, Un esprit sceptique poussé trop loin, Une théorie complotiste comme il en existe sur d'autres sujets (Terre plate, Apollo 11 etc.), Problème d'éducation et/ou d'information"
$ws.Cells.Item(56,11).Value = 8.0
$ws.Cells.Item(56,12).Value = 7.0
$ws.Cells.Item(56,13).Value = 6.0
$ws.Cells.Item(56,14).Value = 6.0
$ws.Cells.Item(56,15).Value = 6.0
$ws.Cells.Item(56,16).Value = "Je partage des références directement à mon entourage (mail, vive voix etc.), Je donne une conférence sur les enjeux énergie-climat (hi hi hi)"
$ws.Cells.Item(56,17).Value = 2.0
$ws.Cells.Item(56,18).Value = "On ne peut pas parler de mouvement religieux dans le sens où même si certains livre font foi dans le mouvement, on ne suit aucun écrit à la règle sans réflexions. "
$ws.Cells.Item(56,19).Value = "Une femme"
$ws.Cells.Item(56,20).Value = 20.0
$ws.Cells.Item(56,21).Value = "En banlieue d'un grand centre urbain"
$ws.Cells.Item(56,22).Value = "Études supérieures courtes (DUT BTS ou licence pro en France, Bachelor à l'étranger)"
$ws.Cells.Item(56,23).Value = "Sciences naturelles / Sciences de l'environnement"
$ws.Cells.Item(56,24).Value = "Plutôt écolo (vélo, transport en commun, limitation de la consommation et notamment de la viande)"
$ws.Cells.Item(56,25).Value = "Local, Bio majoritairement (+ de 50% de ce que tu manges chez toi)"
$ws.Cells.Item(56,26).Value = "Gauche modéree (Parti socialiste ou Génération-s en France)"
$ws.Cells.Item(56,28).Value = "Académicien"

# Row 57 values
$ws.Cells.Item(57,1).Value = 43543.80020091435
$ws.Cells.Item(57,2).Value = "Entre 2 et 5 ans"
$ws.Cells.Item(57,3).Value = "Conférences (Jancovici, Bihouix etc.), Articles de vulgarisation & blogs, Vidéos Youtube de vulgarisation, Publications scientifiques"
$ws.Cells.Item(57,4).Value = 6.0
$ws.Cells.Item(57,5).Value = 10.0
$ws.Cells.Item(57,6).Value = "Va commencer dans les 10 ans qui viennent"
$ws.Cells.Item(57,7).Value = "Très rapide (de l'ordre de 3 à 5 ans)"
$ws.Cells.Item(57,8).Value = "Désabusement"
$ws.Cells.Item(57,9).Value = "Je comprend qu'il y ait des personnes climatosceptiques au sein de la population, Je suis prêt à baisser mon niveau de vie même si cette baisse ne s'opère pas pour les autres"
$ws.Cells.Item(57,10).Value = "Une réaction de personnes ne voulant pas toucher à leur mode de vie, Une stupidité"
$ws.Cells.Item(57,11).Value = 8.0
$ws.Cells.Item(57,12).Value = 7.0
$ws.Cells.Item(57,13).Value = 7.0
$ws.Cells.Item(57,14).Value = 6.0
$ws.Cells.Item(57,15).Value = 6.0
$ws.Cells.Item(57,16).Value = "Je partage des liens sur les réseaux sociaux, Je partage des références directement à mon entourage (mail, vive voix etc.), Je donne une conférence sur les enjeux énergie-climat (hi hi hi)"
$ws.Cells.Item(57,17).Value = 4.0
$ws.Cells.Item(57,19).Value = "Une femme"
$ws.Cells.Item(57,20).Value = 24.0
$ws.Cells.Item(57,21).Value = "En banlieue d'un grand centre urbain"
$ws.Cells.Item(57,22).Value = "Études supérieures longues (ingénieur, école de commerce, Master MBA graduate à l'étranger)"
$ws.Cells.Item(57,23).Value = "Ingénierie industrielle / énergétique"
$ws.Cells.Item(57,24).Value = "Plutôt écolo (vélo, transport en commun, limitation de la consommation et notamment de la viande)"
$ws.Cells.Item(57,25).Value = "Local, Flexitarien"
$ws.Cells.Item(57,26).Value = "Aucun de ces partis, je ne crois pas à la politique"
$ws.Cells.Item(57,28).Value = "Académicien"

# Row 58 values
$ws.Cells.Item(58,1).Value = 43543.84074148149
$ws.Cells.Item(58,2).Value = "Depuis + de 8 ans"
$ws.Cells.Item(58,3).Value = "Conférences (Jancovici, Bihouix etc.), Livres, Vidéos Youtube de vulgarisation"
$ws.Cells.Item(58,4).Value = 5.0
$ws.Cells.Item(58,5).Value = 9.0
$ws.Cells.Item(58,6).Value = "A déjà commencé"
$ws.Cells.Item(58,7).Value = "Plutôt rapide (5 à 10 ans)"
$ws.Cells.Item(58,8).Value = "Paix intérieure"
$ws.Cells.Item(58,9).Value = "Je suis prêt à baisser mon niveau de vie même si cette baisse ne s'opère pas pour les autres, Une transition écologique efficace peut se faire dans un cadre démocratique"
$ws.Cells.Item(58,10).Value = "Un esprit sceptique poussé trop loin, Problème d'éducation et/ou d'information"
$ws.Cells.Item(58,11).Value = 7.0
$ws.Cells.Item(58,12).Value = 6.0
$ws.Cells.Item(58,13).Value = 5.0
$ws.Cells.Item(58,14).Value = 2.0
$ws.Cells.Item(58,15).Value = 4.0
$ws.Cells.Item(58,16).Value = "Je partage des références directement à mon entourage (mail, vive voix etc.), Je donne une conférence sur les enjeux énergie-climat (hi hi hi)"
$ws.Cells.Item(58,17).Value = 5.0
$ws.Cells.Item(58,18).Value = "Cela ressemble énormément à un dogme, mais j'y suis sensible malgré cette impression. Je trouve que nous sommes poussés à nous éloigner de nos ressentis par la société. En prenant le temps de regarder en dehors du cadre imposé par les normes sociales, les enjeux climatique me semble être une évidence."
$ws.Cells.Item(58,19).Value = "Une femme"
$ws.Cells.Item(58,20).Value = 21.0
$ws.Cells.Item(58,21).Value = "En ville dans une ville moyenne"
$ws.Cells.Item(58,22).Value = "Études supérieures longues (ingénieur, école de commerce, Master MBA graduate à l'étranger)"
$ws.Cells.Item(58,23).Value = "Mathématiques / Informatique"
$ws.Cells.Item(58,24).Value = "Plutôt écolo (vélo, transport en commun, limitation de la consommation et notamment de la viande)"
$ws.Cells.Item(58,25).Value = "Local, Bio majoritairement (+ de 50% de ce que tu manges chez toi), Flexitarien"
$ws.Cells.Item(58,26).Value = "Très à gauche (France insoumise ou plus à gauche en France)"
$ws.Cells.Item(58,28).Value = "Académicien"

# Row 59 values
$ws.Cells.Item(59,1).Value = 43543.94531111111
$ws.Cells.Item(59,2).Value = "Depuis + de 8 ans"
$ws.Cells.Item(59,3).Value = "Livres, une intuition venue dès la naissance, grande sensibilité pour les animaux et la protection de la nature, sans vraiment lire de documents scientifiques. quelques documentaires grand public (la planète bleue, la planète blanche)"
$ws.Cells.Item(59,4).Value = 4.0
$ws.Cells.Item(59,5).Value = 7.0
$ws.Cells.Item(59,6).Value = "Va commencer dans les 15 à 20 ans qui viennent"
$ws.Cells.Item(59,7).Value = "Un peu plus lent (de l'ordre de 20 à 30 ans)"
$ws.Cells.Item(59,8).Value = "mélange de tristesse et d'excitation de qu'est ce qu'on va faire ensuite"
$ws.Cells.Item(59,9).Value = "Je suis prêt à baisser mon niveau de vie si cette baisse s'opère pour les autres également, Je suis prêt à baisser mon niveau de vie même si cette baisse ne s'opère pas pour les autres, Une transition écologique efficace peut se faire dans un cadre démocratique, Pour agir efficacement, il faut hiérarchiser les problèmes (perte de biodiversité, dérèglement climatique etc.)"
$ws.Cells.Item(59,10).Value = "Une réaction de personnes ne voulant pas toucher à leur mode de vie, Problème d'éducation et/ou d'information, une fuite pour ne pas faire face au problème"
$ws.Cells.Item(59,11).Value = 7.0
$ws.Cells.Item(59,12).Value = 6.0
$ws.Cells.Item(59,13).Value = 7.0
$ws.Cells.Item(59,14).Value = 6.0
$ws.Cells.Item(59,15).Value = 6.0
$ws.Cells.Item(59,16).Value = "Je partage des références directement à mon entourage (mail, vive voix etc.), Je donne une conférence sur les enjeux énergie-climat (hi hi hi)"
$ws.Cells.Item(59,17).Value = 2.0
$ws.Cells.Item(59,18).Value = "Le militantisme écologique ne cherche pas à rassembler des personnes dans une foi en un dieu mais de faire prendre conscience d'une situation tangible grave et d'un comportement global non adapté. Je pense néanmoins qu'il peut faire appel à des notions communes à la religion comme le respect du vivant."
$ws.Cells.Item(59,19).Value = "Une femme"
$ws.Cells.Item(59,20).Value = 25.0
$ws.Cells.Item(59,21).Value = "En ville dans une grande agglomération"
$ws.Cells.Item(59,22).Value = "Études supérieures longues (ingénieur, école de commerce, Master MBA graduate à l'étranger)"
$ws.Cells.Item(59,23).Value = "Art / Musique / Design / Arts du spectacles"
$ws.Cells.Item(59,24).Value = "Très frugal (flexitarien ou végétalien, AMAP, déplacement doux)"
$ws.Cells.Item(59,25).Value = "Local, Bio majoritairement (+ de 50% de ce que tu manges chez toi), Végétarien"
$ws.Cells.Item(59,26).Value = "Gauche modéree (Parti socialiste ou Génération-s en France)"
$ws.Cells.Item(59,28).Value = "Académicien, investie dans le Pôle culturel"

# Row 60 values
$ws.Cells.Item(60,1).Value = 43544.06678005787
$ws.Cells.Item(60,2).Value = "Entre 5 et 8 ans"
$ws.Cells.Item(60,3).Value = "Conférences (Jancovici, Bihouix etc.), Articles de vulgarisation & blogs, Livres, Vidéos Youtube de vulgarisation, Articles de presse, Publications scientifiques"
$ws.Cells.Item(60,4).Value = 9.0
$ws.Cells.Item(60,5).Value = 7.0
$ws.Cells.Item(60,6).Value = "A déjà commencé"
$ws.Cells.Item(60,7).Value = "Un peu plus lent (de l'ordre de 20 à 30 ans)"
$ws.Cells.Item(60,8).Value = "Paix intérieure"
$ws.Cells.Item(60,9).Value = "On a un fort potentiel d'action à l'échelle individuelle, Je comprend qu'il y ait des personnes climatosceptiques au sein de la population, Une transition écologique efficace peut se faire dans un cadre démocratique"
$ws.Cells.Item(60,10).Value = "Une opinion différente, mais qui a une probabilité d'être juste, Un terme assez condescendant pour catégoriser ceux qui remettent en question le caractère majoritairement anthropique du changement climatique "
$ws.Cells.Item(60,11).Value = 8.0
$ws.Cells.Item(60,12).Value = 2.0
$ws.Cells.Item(60,13).Value = 6.0
$ws.Cells.Item(60,14).Value = 7.0
$ws.Cells.Item(60,15).Value = 4.0
$ws.Cells.Item(60,16).Value = "Je fais profil bas. Trop en parler, c'est devenir prêcheur, et donc desservir la cause., Je donne une conférence sur les enjeux énergie-climat (hi hi hi), Je montre l'exemple"
$ws.Cells.Item(60,17).Value = 4.0
$ws.Cells.Item(60,18).Value = "Dogme religieux fermé aux opinions contraires et qui les rejette fortement "
$ws.Cells.Item(60,19).Value = "Un homme"
$ws.Cells.Item(60,20).Value = 25.0
$ws.Cells.Item(60,21).Value = "En ville dans une grande agglomération"
$ws.Cells.Item(60,22).Value = "Études supérieures longues (ingénieur, école de commerce, Master MBA graduate à l'étranger)"
$ws.Cells.Item(60,23).Value = "Economie / Gestion"
$ws.Cells.Item(60,24).Value = "Plutôt écolo (vélo, transport en commun, limitation de la consommation et notamment de la viande)"
$ws.Cells.Item(60,25).Value = "Local, Bio majoritairement (+ de 50% de ce que tu manges chez toi), Flexitarien"
$ws.Cells.Item(60,26).Value = "Autre"
$ws.Cells.Item(60,28).Value = "Académicien"

# Row 61 values
$ws.Cells.Item(61,1).Value = 43544.51922133102
$ws.Cells.Item(61,2).Value = "Entre 5 et 8 ans"
$ws.Cells.Item(61,3).Value = "Conférences (Jancovici, Bihouix etc.), Le bouche à oreille (conférences, rencontres, associations...)"
$ws.Cells.Item(61,4).Value = 5.0
$ws.Cells.Item(61,5).Value = 8.0
$ws.Cells.Item(61,6).Value = "A déjà commencé"
$ws.Cells.Item(61,7).Value = "Plutôt rapide (5 à 10 ans)"
$ws.Cells.Item(61,8).Value = "Tristesse"
$ws.Cells.Item(61,9).Value = "On a un fort potentiel d'action à l'échelle individuelle, Je suis prêt à baisser mon niveau de vie si cette baisse s'opère pour les autres également, Je suis prêt à baisser mon niveau de vie même si cette baisse ne s'opère pas pour les autres, Une transition écologique efficace peut se faire dans un cadre démocratique, Pour agir efficacement, il faut hiérarchiser les problèmes (perte de biodiversité, dérèglement climatique etc.)"
$ws.Cells.Item(61,10).Value = "Une réaction de personnes ne voulant pas toucher à leur mode de vie, Problème d'éducation et/ou d'information"
$ws.Cells.Item(61,11).Value = 8.0
$ws.Cells.Item(61,12).Value = 7.0
$ws.Cells.Item(61,13).Value = 7.0
$ws.Cells.Item(61,14).Value = 4.0
$ws.Cells.Item(61,15).Value = 4.0
$ws.Cells.Item(61,16).Value = "Je partage des liens sur les réseaux sociaux, Je donne une conférence sur les enjeux énergie-climat (hi hi hi), Je leur fait à manger vegan (et c'est bon niark niark niark)"
$ws.Cells.Item(61,17).Value = 1.0
$ws.Cells.Item(61,18).Value = "Je trouve que les personnes écolos ont des références assez variées et sont touchées par le réchauffement climatique pour différentes causes et à différents dégrés"
$ws.Cells.Item(61,19).Value = "Une femme"
$ws.Cells.Item(61,20).Value = 23.0
$ws.Cells.Item(61,21).Value = "En ville dans une grande agglomération"
$ws.Cells.Item(61,22).Value = "Études supérieures longues (ingénieur, école de commerce, Master MBA graduate à l'étranger)"
$ws.Cells.Item(61,23).Value = "Sciences sociales"
$ws.Cells.Item(61,24).Value = "Plutôt écolo (vélo, transport en commun, limitation de la consommation et notamment de la viande)"
$ws.Cells.Item(61,25).Value = "Bio majoritairement (+ de 50% de ce que tu manges chez toi), Végétalien"
$ws.Cells.Item(61,26).Value = "Autre"
$ws.Cells.Item(61,28).Value = "Académicien"

# --- Copy cell formatting (styles) from row 55 template to new rows ---
# Column A uses the date/time style; other present columns use the standard text/number style.
# Runs are split so we never touch columns absent from the target row (e.g. AA, or a skipped question column).
# Row 56 formatting
$ws.Range("A55:Z55").Copy()
$ws.Range("A56:Z56").PasteSpecial(-4122)
$ws.Range("AB55").Copy()
$ws.Range("AB56").PasteSpecial(-4122)

# Row 57 formatting
$ws.Range("A55:Q55").Copy()
$ws.Range("A57:Q57").PasteSpecial(-4122)
$ws.Range("S55:Z55").Copy()
$ws.Range("S57:Z57").PasteSpecial(-4122)
$ws.Range("AB55").Copy()
$ws.Range("AB57").PasteSpecial(-4122)

# Row 58 formatting
$ws.Range("A55:Z55").Copy()
$ws.Range("A58:Z58").PasteSpecial(-4122)
$ws.Range("AB55").Copy()
$ws.Range("AB58").PasteSpecial(-4122)

# Row 59 formatting
$ws.Range("A55:Z55").Copy()
$ws.Range("A59:Z59").PasteSpecial(-4122)
$ws.Range("AB55").Copy()
$ws.Range("AB59").PasteSpecial(-4122)

# Row 60 formatting
$ws.Range("A55:Z55").Copy()
$ws.Range("A60:Z60").PasteSpecial(-4122)
$ws.Range("AB55").Copy()
$ws.Range("AB60").PasteSpecial(-4122)

# Row 61 formatting
$ws.Range("A55:Z55").Copy()
$ws.Range("A61:Z61").PasteSpecial(-4122)
$ws.Range("AB55").Copy()
$ws.Range("AB61").PasteSpecial(-4122)

$excel.CutCopyMode = 0
Write-Host "Done adding rows 56-61"